# Apply the update described by the diff:
#  - B36: 90658 -> 90792
#  - B37: 90689 -> 90823
#  - Row 38 and Row 39 swap all of their data (columns C..AY), while A and B
#    on each row get their own independent new values.
#  - B40: 88949 -> 89083

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns (as text-like) that Excel would otherwise auto-convert to a
# number/date/time when assigned via .Value2 (e.g. "2023-09-12" -> date
# serial, "1" -> number). These must be forced back to plain text so the
# stored cell type/value matches the original inline string.
$textColumns = @(9, 25, 27)   # I (Antal), Y (Startdatum), AA (Slutdatum)

function Set-CellSmart {
    param($cell, $value, $forceText)

    if ($forceText -and $value -ne $null -and $value -ne "") {
        $cell.NumberFormat = "@"
        $cell.Value2 = $value
        $cell.Style = "Normal"
    } else {
        $cell.Value2 = $value
    }
}

# Capture the existing values of row 38 and row 39 (columns C through AY,
# i.e. column indices 3 through 51) before any writes happen.
$row38Old = @{}
$row39Old = @{}
for ($c = 3; $c -le 51; $c++) {
    $row38Old[$c] = $ws.Cells.Item(38, $c).Value2
    $row39Old[$c] = $ws.Cells.Item(39, $c).Value2
}

# Write row 39's old content into row 38, and row 38's old content into row 39.
for ($c = 3; $c -le 51; $c++) {
    $forceText = $textColumns -contains $c
    Set-CellSmart $ws.Cells.Item(38, $c) $row39Old[$c] $forceText
    Set-CellSmart $ws.Cells.Item(39, $c) $row38Old[$c] $forceText
}

# Independent new Id / Taxonsorteringsordning values for rows 38 and 39.
$ws.Cells.Item(38, 1).Value2 = 112045406
$ws.Cells.Item(38, 2).Value2 = 90816

$ws.Cells.Item(39, 1).Value2 = 112073630
$ws.Cells.Item(39, 2).Value2 = 89058

# Independent Taxonsorteringsordning updates on rows 36, 37 and 40.
$ws.Cells.Item(36, 2).Value2 = 90792
$ws.Cells.Item(37, 2).Value2 = 90823
$ws.Cells.Item(40, 2).Value2 = 89083
